$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15-26: Status (col E) moves from "Error" to "Corrected",
# and Fixed Date (col S) / Closed Date (col T) get filled in with 2011-10-22 (serial 40838).
$ws.Range("E15:E26").Value = "Corrected"

$ws.Range("S15:S26").Value = 40838
# Copy S's format onto T so the Closed Date cells pick up the same date
# number-format / border style as the Fixed Date cells (matches the
# style index Excel produced when the two columns were filled together).
$ws.Range("S15:S26").Copy()
$ws.Range("T15:T26").PasteSpecial(-4122)
$ws.Range("T15:T26").Value = 40838

$excel.CutCopyMode = 0

# Leave the selection where the author left it when they saved the file.
$ws.Range("E16").Select()
